$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data rows for the "Land use" indicators.
# Row 74: id 73, Land use (conventional)
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = "Land use (conventional)"

# Row 75: id 74, Land use (RES)
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = "Land use (RES)"

# Fill the description column after both labels, matching the order in
# which the shared strings were authored.
$ws.Range("C74").Value = "Land use per MWh produced electricity _fossil electricity generation (m2/MWh)"
$ws.Range("C75").Value = "Land use per MWh produced electricity_renewable electricity generation (m2/MWh)"

# Reflect the final cursor/selection position left by the editor.
$ws.Range("B79").Select()
